$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row labels: "<Name>_old" -> "<Name>_FV2410" and
#    "<Name>_new" -> "<Name>_FV2504" (columns A-J and L-U of row 1; column K
#    stays "diff").
# ---------------------------------------------------------------------------
$labels = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $labels[$i] + "_FV2410"
    $ws.Cells.Item(1, $i + 12).Value = $labels[$i] + "_FV2504"
}

# ---------------------------------------------------------------------------
# 2) Turn the used range into an Excel Table ("Table1") covering A1:U66.
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:U66")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1, $null)
$lo.Name = "Table1"

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split/freeze after row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
